$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newUrl = "https://m.media-amazon.com/images/S/compressed.photo.goodreads.com/books/1670363463i/58416952.jpg"
$target = $ws.Range("C11")
$target.Clear()
Write-Host ("count after clear: " + $ws.Hyperlinks.Count)
$ws.Hyperlinks.Add($target, $newUrl, "", $newUrl, $newUrl)
Write-Host ("count after add: " + $ws.Hyperlinks.Count)
